# Generate Report for Handback
# Updates status/timestamp cells across the Overview, zh-cn and de-de sheets
# to reflect the latest handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 795f56e8... row and the 9a3e5ccf... row.
$wsOverview.Range("G4").Value = "2016-09-01 00:18:11"
$wsOverview.Range("G5").Value = "2016-09-01 00:18:11"

# zh-cn sheet: Priority column (E) status changes from "ht" to "mt",
# and the handoff/handback datetimes (H/K) move forward.
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-09-01 00:18:00"
$wsZhCn.Range("H5").Value = "2016-09-01 00:18:00"
$wsZhCn.Range("K4").Value = "2016-09-01 00:18:29"
$wsZhCn.Range("K5").Value = "2016-09-01 00:18:29"

# de-de sheet: handoff datetime (H) and handback datetime (K) move forward.
$wsDeDe.Range("H4").Value = "2016-09-01 00:18:11"
$wsDeDe.Range("H5").Value = "2016-09-01 00:18:11"
$wsDeDe.Range("K4").Value = "2016-09-01 00:18:36"
$wsDeDe.Range("K5").Value = "2016-09-01 00:18:36"
